$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44498
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = 100112045
$ws.Range("G41").Value = "Zapallo"
$ws.Range("H41").Value = "Paine"
$ws.Range("I41").Value = "1a (guarda)"
$ws.Range("J41").Value = 4000
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 100
$ws.Range("M41").Value = 100
$ws.Range("N41").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 100
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
